$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (matches original inline-string formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "63.207.88"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "2.465.02"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "578.75"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "147.15"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").Value = "2.463.56"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "28.46"
$ws.Range("E14").Value = "  +7.66%  "
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").Value = "2.898.83"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").Value = "63.093.88"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "2.458.22"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "331.21"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").Value = "  +11.38%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "66.48"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "1.21"
$ws.Range("E26").Value = "  +20.85%  "
$ws.Range("D27").Value = "650.65"
$ws.Range("E27").Value = "  +9.47%  "
$ws.Range("D28").Value = "8.64"
$ws.Range("E28").Value = "  +4.96%  "
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  +7.48%  "
$ws.Range("D30").Value = "2.573.30"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "8.23"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  +6.19%  "
$ws.Range("D33").Value = "1.89"
$ws.Range("E33").Value = "  +4.48%  "
$ws.Range("D34").Value = "0.139"
$ws.Range("E34").Value = "  +4.99%  "
$ws.Range("D35").Value = "0.0₆0396"
$ws.Range("E35").Value = "  +40.13%  "
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +4.17%  "
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("D40").Value = "0.376"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "18.89"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "152.53"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  +10.60%  "
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").Value = "42.67"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +27.36%  "
$ws.Range("D48").Value = "147.39"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").Value = "20.81"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").Value = "0.610"
$ws.Range("E51").Value = "  +3.01%  "
